$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.550.25'
$ws.Range('D3').Value = '2.253.41'
$ws.Range('E3').Value = '  +4.38%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '232.16'
$ws.Range('E5').Value = '  +1.86%  '
$ws.Range('E6').Value = '  +2.21%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '64.20'
$ws.Range('E7').Value = '  +0.40%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.410'
$ws.Range('E9').Value = '  +3.69%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '59.27'
$ws.Range('E10').Value = '  +2.28%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0900'
$ws.Range('E11').Value = '  +5.21%  '
$ws.Range('E12').Value = '  +1.16%  '
$ws.Range('D13').Value = '2.591.57'
$ws.Range('E13').Value = '  +4.62%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '16.31'
$ws.Range('E14').Value = '  +1.69%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '22.58'
$ws.Range('E15').Value = '  +2.08%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.833'
$ws.Range('E16').Value = '  +2.66%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.69'
$ws.Range('E17').Value = '  +3.10%  '
$ws.Range('D18').Value = '2.253.01'
$ws.Range('E18').Value = '  +5.22%  '
$ws.Range('D19').Value = '41.429.84'
$ws.Range('E19').Value = '  +5.18%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '73.99'
$ws.Range('E20').Value = '  +3.00%  '
$ws.Range('D21').Value = '0.0₃0914'
$ws.Range('E21').Value = '  +7.35%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.20'
$ws.Range('E22').Value = '  +1.60%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '251.82'
$ws.Range('E23').Value = '  +9.56%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.42'
$ws.Range('E25').Value = '  +6.00%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.40'
$ws.Range('E26').Value = '  +1.80%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.84'
$ws.Range('E27').Value = '  +1.82%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '173.58'
$ws.Range('E28').Value = '  +0.66%  '
$ws.Range('E29').Value = '  +2.48%  '
$ws.Range('E30').Value = '  +3.20%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.46'
$ws.Range('E31').Value = '  +3.10%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.80'
$ws.Range('E32').Value = '  +7.93%  '
$ws.Range('E33').Value = '  +2.70%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.03'
$ws.Range('E34').Value = '  +6.31%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.75'
$ws.Range('E35').Value = '  +3.38%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0635'
$ws.Range('E36').Value = '  +3.14%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '7.03'
$ws.Range('E37').Value = '  -1.40%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.84'
$ws.Range('E38').Value = '  +8.03%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.000267'
$ws.Range('E40').Value = '  +70.85%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.15%  '
$ws.Range('B42').Value = 'FTXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.87'
$ws.Range('E42').Value = '  +12.83%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0240'
$ws.Range('E43').Value = '  +4.23%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.85'
$ws.Range('E44').Value = '  +13.28%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '102.78'
$ws.Range('E45').Value = '  -0.53%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '17.86'
$ws.Range('E46').Value = '  +1.54%  '
$ws.Range('E47').Value = '  +4.44%  '
$ws.Range('D48').Value = '1.510.75'
$ws.Range('E48').Value = '  -1.27%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0944'
$ws.Range('E49').Value = '  +1.47%  '
$ws.Range('E50').Value = '  +2.37%  '
$ws.Range('E51').Value = '  -0.83%  '
